# Insert a new weekly price record before the existing row 119,
# shifting all subsequent rows (119-228) down to (120-229).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows(119).Insert()

$ws.Range("A119").Value = 5
$ws.Range("B119").Value = "Macroferia Regional de Talca"
$ws.Range("C119").Value = "Maule"
$ws.Range("D119").Value = 44512
$ws.Range("E119").Value = 7
$ws.Range("F119").Value = 100112032
$ws.Range("G119").Value = "Zapallo italiano"
$ws.Range("H119").Value = "Sin especificar"
$ws.Range("I119").Value = "Primera"
$ws.Range("J119").Value = 500
$ws.Range("K119").Value = 6000
$ws.Range("L119").Value = 6000
$ws.Range("M119").Value = 6000
$ws.Range("N119").Value = "`$/caja 60 unidades"
$ws.Range("O119").Value = "Región del Maule"
$ws.Range("P119").Value = 100
$ws.Range("Q119").Value = 60
$ws.Range("R119").Value = "Hortaliza"
